$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric cell values (rows 2-9, financial data corrections)

# Row 2
$ws.Range("D2").Value = 506
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 93
$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 76
$ws.Range("K2").Value = 1155
$ws.Range("L2").Value = 153
$ws.Range("M2").Value = 1002
$ws.Range("N2").Value = 1002
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = -110
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = -37
$ws.Range("T2").Value = 27
$ws.Range("U2").Value = -137
$ws.Range("W2").Value = 15.83
$ws.Range("X2").Value = 15.03
$ws.Range("Y2").Value = 7.74
$ws.Range("Z2").Value = 6.27
$ws.Range("AA2").Value = 15.31
$ws.Range("AB2").Value = 2559.7
$ws.Range("AC2").Value = 784
$ws.Range("AD2").Value = 15.68
$ws.Range("AE2").Value = 17800
$ws.Range("AF2").Value = 0.6899999999999999
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 4.88
$ws.Range("AI2").Value = 44.36
$ws.Range("AJ2").Value = 9702706

# Row 3
$ws.Range("D3").Value = 594
$ws.Range("E3").Value = 84
$ws.Range("F3").Value = 84
$ws.Range("G3").Value = 91
$ws.Range("H3").Value = 78
$ws.Range("I3").Value = 78
$ws.Range("K3").Value = 1196
$ws.Range("L3").Value = 151
$ws.Range("M3").Value = 1045
$ws.Range("N3").Value = 1045
$ws.Range("P3").Value = 50
$ws.Range("Q3").Value = 163
$ws.Range("R3").Value = -91
$ws.Range("S3").Value = -34
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = 156
$ws.Range("W3").Value = 14.16
$ws.Range("X3").Value = 13.11
$ws.Range("Y3").Value = 7.61
$ws.Range("Z3").Value = 6.62
$ws.Range("AA3").Value = 14.42
$ws.Range("AB3").Value = 2648.86
$ws.Range("AC3").Value = 802
$ws.Range("AD3").Value = 17.51
$ws.Range("AE3").Value = 18576
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 660
$ws.Range("AH3").Value = 4.7
$ws.Range("AI3").Value = 47.69
$ws.Range("AJ3").Value = 9702706

# Row 4
$ws.Range("D4").Value = 409
$ws.Range("E4").Value = 46
$ws.Range("F4").Value = 46
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 45
$ws.Range("K4").Value = 1160
$ws.Range("L4").Value = 107
$ws.Range("M4").Value = 1053
$ws.Range("N4").Value = 1053
$ws.Range("P4").Value = 50
$ws.Range("Q4").Value = 46
$ws.Range("R4").Value = 148
$ws.Range("S4").Value = -37
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = 41
$ws.Range("W4").Value = 11.3
$ws.Range("X4").Value = 11.01
$ws.Range("Y4").Value = 4.29
$ws.Range("Z4").Value = 3.82
$ws.Range("AA4").Value = 10.19
$ws.Range("AB4").Value = 2662.82
$ws.Range("AC4").Value = 464
$ws.Range("AD4").Value = 24.78
$ws.Range("AE4").Value = 18710
$ws.Range("AF4").Value = 0.61
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 4.35
$ws.Range("AI4").Value = 62.47
$ws.Range("AJ4").Value = 9702706

# Row 5
$ws.Range("D5").Value = 414
$ws.Range("E5").Value = 52
$ws.Range("F5").Value = 52
$ws.Range("G5").Value = 59
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 53
$ws.Range("K5").Value = 1192
$ws.Range("L5").Value = 113
$ws.Range("M5").Value = 1078
$ws.Range("N5").Value = 1078
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = 98
$ws.Range("R5").Value = -136
$ws.Range("S5").Value = -28
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 94
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 12.58
$ws.Range("X5").Value = 12.8
$ws.Range("Y5").Value = 4.97
$ws.Range("Z5").Value = 4.51
$ws.Range("AA5").Value = 10.49
$ws.Range("AB5").Value = 2716.5
$ws.Range("AC5").Value = 546
$ws.Range("AD5").Value = 22.06
$ws.Range("AE5").Value = 19166
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 580
$ws.Range("AH5").Value = 4.81
$ws.Range("AI5").Value = 61.57
$ws.Range("AJ5").Value = 9702706

# Row 6
$ws.Range("D6").Value = 421
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = 60
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = 63
$ws.Range("I6").Value = 63
$ws.Range("K6").Value = 1258
$ws.Range("L6").Value = 151
$ws.Range("M6").Value = 1107
$ws.Range("N6").Value = 1107
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 83
$ws.Range("R6").Value = -136
$ws.Range("S6").Value = -33
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 78
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 14.31
$ws.Range("X6").Value = 15.05
$ws.Range("Y6").Value = 5.8
$ws.Range("Z6").Value = 5.18
$ws.Range("AA6").Value = 13.63
$ws.Range("AB6").Value = 2774.02
$ws.Range("AC6").Value = 654
$ws.Range("AD6").Value = 19.5
$ws.Range("AE6").Value = 19680
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 650
$ws.Range("AH6").Value = 5.1
$ws.Range("AI6").Value = 57.65
$ws.Range("AJ6").Value = 9702706

# Row 7
$ws.Range("D7").Value = 398
$ws.Range("E7").Value = 35
$ws.Range("G7").Value = 54
$ws.Range("H7").Value = 46
$ws.Range("I7").Value = 46
$ws.Range("K7").Value = 1249
$ws.Range("L7").Value = 132
$ws.Range("M7").Value = 1117
$ws.Range("N7").Value = 1117
$ws.Range("P7").Value = 50
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = -17
$ws.Range("S7").Value = -37
$ws.Range("T7").Value = 5
$ws.Range("W7").Value = 8.789999999999999
$ws.Range("X7").Value = 11.56
$ws.Range("Y7").Value = 4.14
$ws.Range("Z7").Value = 3.67
$ws.Range("AA7").Value = 11.82
$ws.Range("AC7").Value = 474
$ws.Range("AD7").Value = 25.94
$ws.Range("AE7").Value = 19852
$ws.Range("AF7").Value = 0.62
$ws.Range("AG7").Value = 600
$ws.Range("AH7").Value = 4.88
$ws.Range("AI7").Value = 126.56

# Row 8
$ws.Range("D8").Value = 450
$ws.Range("E8").Value = 42
$ws.Range("G8").Value = 60
$ws.Range("H8").Value = 51
$ws.Range("I8").Value = 51
$ws.Range("K8").Value = 1277
$ws.Range("L8").Value = 143
$ws.Range("M8").Value = 1134
$ws.Range("N8").Value = 1134
$ws.Range("P8").Value = 50
$ws.Range("Q8").Value = 64
$ws.Range("R8").Value = -19
$ws.Range("S8").Value = -34
$ws.Range("T8").Value = 6
$ws.Range("W8").Value = 9.33
$ws.Range("X8").Value = 11.33
$ws.Range("Y8").Value = 4.53
$ws.Range("Z8").Value = 4.04
$ws.Range("AA8").Value = 12.61
$ws.Range("AC8").Value = 526
$ws.Range("AD8").Value = 23.4
$ws.Range("AE8").Value = 20154
$ws.Range("AF8").Value = 0.61
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 4.88
$ws.Range("AI8").Value = 114.15

# Row 9
$ws.Range("D9").Value = 489
$ws.Range("E9").Value = 47
$ws.Range("G9").Value = 66
$ws.Range("H9").Value = 56
$ws.Range("I9").Value = 56
$ws.Range("K9").Value = 1308
$ws.Range("L9").Value = 151
$ws.Range("M9").Value = 1156
$ws.Range("N9").Value = 1156
$ws.Range("P9").Value = 50
$ws.Range("Q9").Value = 69
$ws.Range("R9").Value = -19
$ws.Range("S9").Value = -34
$ws.Range("T9").Value = 6
$ws.Range("W9").Value = 9.609999999999999
$ws.Range("X9").Value = 11.45
$ws.Range("Y9").Value = 4.89
$ws.Range("Z9").Value = 4.33
$ws.Range("AA9").Value = 13.06
$ws.Range("AC9").Value = 577
$ws.Range("AD9").Value = 21.31
$ws.Range("AE9").Value = 20545
$ws.Range("AF9").Value = 0.6
$ws.Range("AG9").Value = 600
$ws.Range("AH9").Value = 4.88
$ws.Range("AI9").Value = 103.96

# Clear cells removed entirely in this row (columns no longer reported)

# Row 2
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("V2").ClearContents()

# Row 3
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 7
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("U9").ClearContents()
